# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for the first tracked
# file (66a28a16-a372-45bd-affb-1fc4add6ce51.md) as part of a newer
# handback-status report generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-18 18:50:45"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-18 18:50:39"
$zhcn.Range("K2").Value = "2016-08-18 18:50:57"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-18 18:50:45"
$dede.Range("K2").Value = "2016-08-18 18:51:13"
